# "changed currency statement format"
#
# The foreign-currency columns (Döviz Borç / Döviz Alacak / Döviz Bakiye)
# are removed from the statement header row (F4:H4) - the cells are
# cleared but keep their existing formatting/borders. Since those three
# shared strings are no longer referenced anywhere, they drop out of the
# workbook's shared string table automatically on save.
#
# The totals cells next to the Vergi No (M2:N2) are re-aligned from the
# previous top alignment to a vertically centered alignment.
#
# Finally, the sheet's remembered selection/active cell is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Döviz Borç / Döviz Alacak / Döviz Bakiye" header labels
# (formatting/styles on F4:H4 stay untouched).
$ws.Range("F4:H4").ClearContents()

# Vertically center M2:N2 (previously vertically top-aligned).
$ws.Range("M2:N2").VerticalAlignment = -4108  # xlCenter

# Update the stored selection / active cell for the sheet.
$ws.Range("H14").Select()
